$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = "PyTorch 심층신경망(DNN) 모델 생성 후 Fashion MNIST 이미지 분류기 생성, 학습, 예측, 검증 성능 측정하기"
$ws.Range("E4").Value = "https://teddylee777.github.io/pytorch/pytorch-dnn-fashion-mnist"

$ws.Range("D6").Value = "[Object Detection] YOLO v5, v6 Loss"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Object-Detection-YOLO-v5-v6-Loss"
